# Generate Report for Handback
# Update the "Xliff Generate Date" / handoff / handback timestamp cells
# to reflect the newer report generation run.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: Latest HO Xliff Generate Date for be4ba6f3-...md
$wsOverview.Range("G2").Value = "2016-09-06 07:16:52"

# zh-cn sheet: Correspond Handoff Datetime / Correspond Handback DateTime
# for be4ba6f3-...md
$wsZhCn.Range("H2").Value = "2016-09-06 07:16:47"
$wsZhCn.Range("K2").Value = "2016-09-06 07:17:09"

# de-de sheet: Correspond Handoff Datetime (shares value with Overview!G2)
# and Correspond Handback DateTime for be4ba6f3-...md
$wsDeDe.Range("H2").Value = "2016-09-06 07:16:52"
$wsDeDe.Range("K2").Value = "2016-09-06 07:17:18"
